$d = $word.ActiveDocument

function Get-ParagraphForBookmark($bookmarkName) {
    $b = $d.Bookmarks.Item($bookmarkName)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($b.Start -ge $p.Range.Start -and $b.Start -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

# 1) Strip the dynamic value runs that follow each of these bookmarks,
#    leaving the bookmarkStart/bookmarkEnd pair but no trailing run text.
$namesToClear = @(
    "localidad",
    "elDiaEstara",
    "grados",
    "precipitaciones",
    "humedad",
    "viento",
    "eventos",
    "emisor",
    "asunto",
    "descripcion",
    "priact",
    "segunact",
    "terceract"
)

foreach ($name in $namesToClear) {
    $b = $d.Bookmarks.Item($name)
    $p = Get-ParagraphForBookmark $name
    $delStart = $b.End
    $delEnd = $p.Range.End - 1
    if ($delEnd -gt $delStart) {
        $r = $d.Range($delStart, $delEnd)
        $r.Delete()
    }
}

# 2) Turn the empty paragraph right after "Tus correos" into the
#    "no hay correo" placeholder paragraph: underlined paragraph mark +
#    a brand new "nohaycorreo" bookmark.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Tus correos*") {
        $target = $d.Paragraphs.Item($i + 1)

        # Insert a temporary character so the run/paragraph-mark formatting
        # can be applied, then add the bookmark, then remove the temp char.
        $ins = $d.Range($target.Range.Start, $target.Range.Start)
        $ins.InsertAfter("X")

        $p2 = $d.Paragraphs.Item($i + 1)
        $p2.Range.Font.Underline = 1

        $p3 = $d.Paragraphs.Item($i + 1)
        $bmRange = $d.Range($p3.Range.Start, $p3.Range.Start)
        $d.Bookmarks.Add("nohaycorreo", $bmRange)

        $p4 = $d.Paragraphs.Item($i + 1)
        $delRange = $d.Range($p4.Range.Start, $p4.Range.End - 1)
        $delRange.Delete()
        break
    }
}
